{"js": "// MRD-1896 Only show PPCS query emails in one question\n//\n// Two merge-field placeholders in the document currently read\n// \"{{ppcs_query_emails}}\" (in the text cached as the result of a Word\n// FORMTEXT field). They need distinct, question-specific field names so\n// the same placeholder isn't shown for both questions:\n//   - the one under \"25. Probation details\" (preceded by the\n//     \"Fax Number:\" field)            -> {{completed_by_ppcs_query_emails}}\n//   - the one under \"Supervising practitioner\" details (preceded by the\n//     \"LDU: {{supervising_practitioner_local_delivery_unit}}\" field)\n//                                      -> {{supervising_practitioner_ppcs_query_emails}}\n//\n// We locate each occurrence with a search, disambiguate using the\n// paragraph two steps back (the nearest preceding paragraph with\n// distinguishing text), and insert the new prefix right after the\n// opening \"{{\" so the trailing \"}}\" + spaces are left untouched.\n\nconst searchText = \"{{ppcs_query_emails}\";\n\nconst results = context.document.body.search(searchText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length !== 2) {\n  throw new Error(\"Expected 2 occurrences of \" + searchText + \", found \" + results.items.length);\n}\n\n// For each match, walk back to the nearest ancestor paragraph that holds\n// distinguishing context text, so we know which prefix belongs where.\nconst contextRanges = [];\nfor (let i = 0; i < results.items.length; i++) {\n  const para = results.items[i].paragraphs.getFirst();\n  const prev1 = para.getPreviousOrNullObject();\n  prev1.load(\"isNullObject\");\n  contextRanges.push({ match: results.items[i], prev1 });\n}\nawait context.sync();\n\nconst prev2List = [];\nfor (const { prev1 } of contextRanges) {\n  const prev2 = prev1.isNullObject ? null : prev1.getPreviousOrNullObject();\n  if (prev2) {\n    prev2.load(\"text,isNullObject\");\n  }\n  prev2List.push(prev2);\n}\nawait context.sync();\n\nfor (let i = 0; i < contextRanges.length; i++) {\n  const { match } = contextRanges[i];\n  const prev2 = prev2List[i];\n  const nearbyText = prev2 && !prev2.isNullObject ? prev2.text : \"\";\n\n  let prefix;\n  if (nearbyText.indexOf(\"supervising_practitioner_local_delivery_unit\") !== -1) {\n    prefix = \"supervising_practitioner_\";\n  } else {\n    prefix = \"completed_by_\";\n  }\n\n  // Insert the new prefix right after the literal \"{{\" so the existing\n  // closing brace(s) and trailing spaces (which live in the following\n  // run) are left exactly as they were.\n  const openBrace = match.search(\"{{\", { matchCase: true }).getFirstOrNullObject();\n  openBrace.load(\"isNullObject\");\n  await context.sync();\n\n  const insertionPoint = openBrace.getRange(\"End\");\n  insertionPoint.insertText(prefix, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# MRD-1896 Only show PPCS query emails in one question\n#\n# Two merge-field placeholders currently read \"{{ppcs_query_emails}}\"\n# (the cached display text of a Word FORMTEXT field). Give each one a\n# distinct, question-specific field name so the same placeholder isn't\n# shown for both questions:\n#   - the one under \"25. Probation details\" (preceded by the\n#     \"Fax Number:\" field)            -> {{completed_by_ppcs_query_emails}}\n#   - the one under \"Supervising practitioner\" details (preceded by the\n#     \"LDU: {{supervising_practitioner_local_delivery_unit}}\" field)\n#                                      -> {{supervising_practitioner_ppcs_query_emails}}\n#\n# Each paragraph is located, disambiguated using the nearest preceding\n# paragraph with distinguishing text, and the new prefix is inserted\n# right after the opening \"{{\" so the existing \"}}\" and trailing spaces\n# are left untouched.\n\n$d = $word.ActiveDocument\n\n$paragraphs = @()\nforeach ($p in $d.Paragraphs) {\n  $paragraphs += $p\n}\n\n$searchText = \"{{ppcs_query_emails}\"\n$count = $paragraphs.Count\n$matchCount = 0\n\nfor ($i = 0; $i -lt $count; $i++) {\n  $paraText = $paragraphs[$i].Range.Text\n\n  if ($paraText -like \"*$searchText*\") {\n    $matchCount = $matchCount + 1\n\n    $nearbyText = \"\"\n    if ($i -ge 2) {\n      $nearbyText = $paragraphs[$i - 2].Range.Text\n    }\n\n    if ($nearbyText -like \"*supervising_practitioner_local_delivery_unit*\") {\n      $prefix = \"supervising_practitioner_\"\n    } else {\n      $prefix = \"completed_by_\"\n    }\n\n    $target = $paragraphs[$i].Range\n    $found = $target.Find.Execute($searchText)\n\n    if ($found) {\n      $insertRange = $d.Range($target.Start + 2, $target.Start + 2)\n      $insertRange.InsertAfter($prefix)\n    }\n  }\n}\n\nWrite-Output (\"matches=\" + $matchCount)\n"}
